$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain plain text even for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.997.92'
$ws.Range("E2").Value = '  +0.52%  '

# Row 3
$ws.Range("D3").Value = '1.883.64'
$ws.Range("E3").Value = '  -0.23%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("E5").Value = '  -2.33%  '

# Row 6
$ws.Range("D6").Value = '242.16'
$ws.Range("E6").Value = '  -0.06%  '

# Row 7
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").Value = '0.3158'
$ws.Range("E8").Value = '  +1.05%  '

# Row 9
$ws.Range("D9").Value = '0.07161'
$ws.Range("E9").Value = '  +0.62%  '

# Row 10
$ws.Range("D10").Value = '24.61'
$ws.Range("E10").Value = '  -2.77%  '

# Row 11
$ws.Range("D11").Value = '0.08327'
$ws.Range("E11").Value = '  -1.95%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7539'
$ws.Range("E12").Value = '  -0.75%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.906.98'
$ws.Range("E13").Value = '  -0.76%  '

# Row 14
$ws.Range("D14").Value = '5.390'
$ws.Range("E14").Value = '  +0.37%  '

# Row 15
$ws.Range("D15").Value = '92.52'
$ws.Range("E15").Value = '  -0.89%  '

# Row 16
$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").Value = '6.130'
$ws.Range("E16").Value = '  +0.01%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.990.46'
$ws.Range("E17").Value = '  +0.34%  '

# Row 18
$ws.Range("D18").Value = '249.33'
$ws.Range("E18").Value = '  +2.50%  '

# Row 19
$ws.Range("D19").Value = '13.51'
$ws.Range("E19").Value = '  -1.45%  '

# Row 20
$ws.Range("D20").Value = '0.000007840'
$ws.Range("E20").Value = '  +0.18%  '

# Row 21
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.175.80'
$ws.Range("E21").Value = '  +1.70%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.12%  '

# Row 23
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.12%  '

# Row 24
$ws.Range("D24").Value = '7.879'
$ws.Range("E24").Value = '  -1.55%  '

# Row 25
$ws.Range("E25").Value = '  -1.77%  '

# Row 26
$ws.Range("D26").Value = '9.254'
$ws.Range("E26").Value = '  -1.24%  '

# Row 27
$ws.Range("D27").Value = '164.24'
$ws.Range("E27").Value = '  +0.72%  '

# Row 28
$ws.Range("D28").Value = '18.62'
$ws.Range("E28").Value = '  -0.48%  '

# Row 29
$ws.Range("D29").Value = '2.043'
$ws.Range("E29").Value = '  +0.59%  '

# Row 30
$ws.Range("D30").Value = '1.474'
$ws.Range("E30").Value = '  -0.56%  '

# Row 31
$ws.Range("D31").Value = '4.552'
$ws.Range("E31").Value = '  +0.93%  '

# Row 32
$ws.Range("D32").Value = '1.533'
$ws.Range("E32").Value = '  +0.06%  '

# Row 33
$ws.Range("D33").Value = '4.174'
$ws.Range("E33").Value = '  +0.27%  '

# Row 34
$ws.Range("D34").Value = '0.05320'
$ws.Range("E34").Value = '  -1.85%  '

# Row 35
$ws.Range("D35").Value = '1.246'
$ws.Range("E35").Value = '  +0.41%  '

# Row 36
$ws.Range("D36").Value = '0.7664'
$ws.Range("E36").Value = '  +1.94%  '

# Row 37
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  -0.28%  '

# Row 38
$ws.Range("D38").Value = '2.724'
$ws.Range("E38").Value = '  +0.48%  '

# Row 39
$ws.Range("D39").Value = '0.01951'
$ws.Range("E39").Value = '  +0.20%  '

# Row 40
$ws.Range("E40").Value = '  -0.41%  '

# Row 41
$ws.Range("D41").Value = '0.4540'
$ws.Range("E41").Value = '  +1.58%  '

# Row 42
$ws.Range("D42").Value = '1.098.54'
$ws.Range("E42").Value = '  -0.03%  '

# Row 43
$ws.Range("D43").Value = '6.038'
$ws.Range("E43").Value = '  -1.17%  '

# Row 44
$ws.Range("D44").Value = '72.26'
$ws.Range("E44").Value = '  -0.57%  '

# Row 45
$ws.Range("D45").Value = '0.8755'
$ws.Range("E45").Value = '  +1.60%  '

# Row 46
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.21%  '

# Row 47
$ws.Range("D47").Value = '104.14'
$ws.Range("E47").Value = '  +1.70%  '

# Row 48
$ws.Range("D48").Value = '1.848'
$ws.Range("E48").Value = '  -0.57%  '

# Row 49
$ws.Range("D49").Value = '7.536'
$ws.Range("E49").Value = '  -2.40%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.061.97'
$ws.Range("E50").Value = '  +1.11%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.545'
$ws.Range("E51").Value = '  -2.01%  '
